$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new TPM data only has 6 data rows (old rows 2-7 get rewritten below);
# the old trailing rows 8, 9 and 10 no longer exist, so remove them.
# Delete from the bottom up so row numbers of the remaining rows don't shift
# out from under us.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Row 2: FAPs -> Dlk1 -> Notch2 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dlk1"
$ws.Range("C2").Value = "Notch2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.101448
$ws.Range("H2").Value = 6.304344
$ws.Range("I2").Value = 0.5480341737688159
$ws.Range("J2").Value = 0.5480341737688159
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9705896666666667
$ws.Range("N2").Value = 2.911769
$ws.Range("O2").Value = 0.02073452941466921
$ws.Range("P2").Value = 0.02073452941466921
$ws.Range("Q2").Value = 2.039643713837334
$ws.Range("R2").Value = 18.356793424536
$ws.Range("S2").Value = 0.01136323069625345
$ws.Range("T2").Value = 0.01136323069625345

# Row 3: FAPs -> Dlk1 -> Notch2 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Dlk1"
$ws.Range("C3").Value = "Notch2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.101448
$ws.Range("H3").Value = 6.304344
$ws.Range("I3").Value = 0.5480341737688159
$ws.Range("J3").Value = 0.5480341737688159
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 26.34807
$ws.Range("N3").Value = 79.04420999999999
$ws.Range("O3").Value = 0.5628689972673966
$ws.Range("P3").Value = 0.5628689972673966
$ws.Range("Q3").Value = 55.36909900535999
$ws.Range("R3").Value = 498.32189104824
$ws.Range("S3").Value = 0.3084714458575196
$ws.Range("T3").Value = 0.3084714458575196

# Row 4: FAPs -> Dlk1 -> Notch2 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Dlk1"
$ws.Range("C4").Value = "Notch2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.101448
$ws.Range("H4").Value = 6.304344
$ws.Range("I4").Value = 0.5480341737688159
$ws.Range("J4").Value = 0.5480341737688159
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 19.49164633333333
$ws.Range("N4").Value = 58.47493899999999
$ws.Range("O4").Value = 0.4163964733179342
$ws.Range("P4").Value = 0.4163964733179341
$ws.Range("Q4").Value = 40.96068120389067
$ws.Range("R4").Value = 368.646130835016
$ws.Range("S4").Value = 0.2281994972150428
$ws.Range("T4").Value = 0.2281994972150428

# Row 5: MuSCs -> Dlk1 -> Notch2 -> ECs
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Dlk1"
$ws.Range("C5").Value = "Notch2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.733072
$ws.Range("H5").Value = 5.199216
$ws.Range("I5").Value = 0.4519658262311841
$ws.Range("J5").Value = 0.4519658262311841
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.9705896666666667
$ws.Range("N5").Value = 2.911769
$ws.Range("O5").Value = 0.02073452941466921
$ws.Range("P5").Value = 0.02073452941466921
$ws.Range("Q5").Value = 1.682101774789333
$ws.Range("R5").Value = 15.138915973104
$ws.Range("S5").Value = 0.009371298718415761
$ws.Range("T5").Value = 0.009371298718415761

# Row 6: MuSCs -> Dlk1 -> Notch2 -> FAPs
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Dlk1"
$ws.Range("C6").Value = "Notch2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.733072
$ws.Range("H6").Value = 5.199216
$ws.Range("I6").Value = 0.4519658262311841
$ws.Range("J6").Value = 0.4519658262311841
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 26.34807
$ws.Range("N6").Value = 79.04420999999999
$ws.Range("O6").Value = 0.5628689972673966
$ws.Range("P6").Value = 0.5628689972673966
$ws.Range("Q6").Value = 45.66310237103999
$ws.Range("R6").Value = 410.96792133936
$ws.Range("S6").Value = 0.254397551409877
$ws.Range("T6").Value = 0.254397551409877

# Row 7: MuSCs -> Dlk1 -> Notch2 -> MuSCs
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Dlk1"
$ws.Range("C7").Value = "Notch2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.733072
$ws.Range("H7").Value = 5.199216
$ws.Range("I7").Value = 0.4519658262311841
$ws.Range("J7").Value = 0.4519658262311841
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 19.49164633333333
$ws.Range("N7").Value = 58.47493899999999
$ws.Range("O7").Value = 0.4163964733179342
$ws.Range("P7").Value = 0.4163964733179341
$ws.Range("Q7").Value = 33.78042649420266
$ws.Range("R7").Value = 304.0238384478239
$ws.Range("S7").Value = 0.1881969761028913
$ws.Range("T7").Value = 0.1881969761028913
